$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9943628311157227
$ws.Range("B1").Value = 1.299428224563599
$ws.Range("C1").Value = 1.129521012306213
$ws.Range("D1").Value = 1.166590690612793
$ws.Range("E1").Value = 1.27840518951416
